# Apply the "make bound-capped uncertainty distributions triangular; make
# sorbate_c,d and sorbate_g,h prod caps same as rest" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: Seed train fermentation ratio -> Triangular shape -----------
$ws.Range("F21").Value = "Triangular"

# --- Row 31: Centrifuge solids recovery -> Triangular shape --------------
$ws.Range("F31").Value = "Triangular"

# --- Row 38: Hydrogenation TAL-to-HMP conversion --------------------------
# Shape becomes Triangular, bound-capped Lower/Midpoint formulas become a
# simple triangular distribution (0.9*baseline / baseline), and the load
# statement is renamed from the old HMTHP reaction to the HMP reaction.
$ws.Range("F38").Value = "Triangular"
$ws.Range("G38").Formula = "=0.9*E38"
$ws.Range("H38").Formula = "=E38"
$ws.Range("K38").Value = "R401.TAL_to_HMP_rxn.X = x"

# --- Rows 40-45: "Dehydration" renamed to "Etherification & hydrolysis" --
$ws.Range("A40").Value = "Etherification & hydrolysis catalyst Amberlyst70:HMP ratio"
$ws.Range("A41").Value = "Etherification & hydrolysis reaction time"
$ws.Range("A42").Value = "Etherification & hydrolysis temperature"
$ws.Range("A43").Value = "Etherification & hydrolysis HMP-to-PSA conversion"
$ws.Range("K43").Value = "R402.HMP_to_PSA_rxn.X = x"
$ws.Range("A44").Value = "Etherification & hydrolysis pressure"
$ws.Range("A45").Value = "Etherification & hydrolysis spent catalyst Amberlyst70 replacement rate"

# --- Row 47: Ring-opening & hydrolysis PSA-to-KS conversion ---------------
# Shape becomes Triangular, and the Lower/Midpoint bounds become a simple
# triangular distribution tied to the baseline (same pattern as row 38).
$ws.Range("F47").Value = "Triangular"
$ws.Range("G47").Formula = "=0.9*E47"
$ws.Range("H47").Formula = "=E47"

# --- View state: scroll down and select F47:H47 (matches author's diff) --
$ws.Range("F47:H47").Select()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
